$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:D4").Value = 0
$ws.Range("B5:C5").Value = 0
$ws.Range("D5").Value = -0.7134594315441092
$ws.Range("B6:D8").Value = 0
$ws.Range("B9:C9").Value = 0
$ws.Range("D9").Value = -0.6775800748059198
